# Chapter10.3.2/position_data.xlsx — refresh the post/date-stamped values
# in the shared strings used by Sheet1 (A1: "Post..." id, A2: "PostCode..." id),
# bumping the embedded date from 03/05/2022 to 01/17/2023.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "Post01172023"
$ws.Range("A2").Value = "PostCode01172023"
